$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program list")
if ($null -eq $ws) { $ws = $wb.ActiveSheet }

# Update row 2 with new program data
$ws.Range("B2").Value = "SAP Development Foundation"
$ws.Range("C2").Value = "6AIV3,4F4Q0"
$ws.Range("D2").Value = $false

# Update row 3 with new program data
$ws.Range("B3").Value = "Cloud Programming Foundation"
$ws.Range("C3").Value = "2724D,46NM9,5H07V"
$ws.Range("D3").Value = $false

# Remove the leftover blank row 4
$ws.Rows.Item(4).Delete()

# Column A has no data for these programs - fully clear the cells (content + format)
$ws.Range("A2:A3").Clear()

# Clear the body rows' explicit cell styling (no longer carries the bordered style)
$ws.Range("B2:D3").Style = "Normal"

# Widen column B to fit the new, longer program names
$ws.Columns.Item(2).ColumnWidth = 28.86

# Move the active selection to C3, matching the saved view state
$ws.Range("C3").Select()
